$p = $ppt.ActivePresentation

# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> currently the stock "Office Theme" palette,
#                             used only by the Notes Master.
#   ppt/theme/theme2.xml  -> currently the "Integral" palette, used by the
#                             Slide Master (and therefore by every slide).
#
# The target revision swaps the two palettes: the Slide Master (and the
# rest of the deck) should end up using the plain "Office" colors, while
# the Notes Master ends up with the former "Integral" colors.
#
# The only theme surface the PowerPoint object model exposes here is the
# presentation's active ThemeColorScheme (reached from a Slide), which is
# backed by the Slide Master's theme part (ppt/theme/theme2.xml). We drive
# all twelve color slots to the values the "Office" theme used to hold, in
# slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
